$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 460.2
$ws.Range("J19").Value = 460.2
$ws.Range("L19").Value = 460.2
$ws.Range("N19").Value = -810.2
$ws.Range("H40").Value = 4416
$ws.Range("I40").Value = 4000.25
$ws.Range("K40").Value = 4000.25
$ws.Range("M40").Value = -3825.25
$ws.Range("H51").Value = 250000300
$ws.Range("I51").Value = 250000300
$ws.Range("K51").Value = 250000300
$ws.Range("M51").Value = -249999816
$ws.Range("H100").Value = 1287.3334
$ws.Range("I100").Value = 1220
$ws.Range("J100").Value = 1725
$ws.Range("K100").Value = 1220
$ws.Range("L100").Value = 1725
$ws.Range("M100").Value = -679
$ws.Range("N100").Value = -2807
$ws.Range("H113").Value = 3955.7273
$ws.Range("I113").Value = 3249
$ws.Range("J113").Value = 4359.5713
$ws.Range("K113").Value = 3249
$ws.Range("L113").Value = 4359.5713
$ws.Range("M113").Value = 5
$ws.Range("N113").Value = -10867.5713
$ws.Range("H116").Value = 13999.8
$ws.Range("I116").Value = 2999
$ws.Range("K116").Value = 2999
$ws.Range("M116").Value = 443
$ws.Range("H132").Value = 4029.3684
$ws.Range("I132").Value = 4346.147
$ws.Range("K132").Value = 13038.441
$ws.Range("M132").Value = -10508.441

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1865.3572
$ws.Range("I32").Value = 1886.4634
$ws.Range("K32").Value = 1886.4634
$ws.Range("M32").Value = -1599.4634
$ws.Range("H61").Value = 2324.5789
$ws.Range("I61").Value = 1591.7826
$ws.Range("K61").Value = 1591.7826
$ws.Range("M61").Value = -1379.7826
$ws.Range("H122").Value = 2559.1052
$ws.Range("I122").Value = 2614
$ws.Range("J122").Value = 2266.3333
$ws.Range("K122").Value = 7842
$ws.Range("L122").Value = 6798.999899999999
$ws.Range("M122").Value = -5392
$ws.Range("N122").Value = -11698.9999
$ws.Range("H136").Value = 2324.5789
$ws.Range("I136").Value = 1591.7826
$ws.Range("K136").Value = 4775.3478
$ws.Range("M136").Value = -2225.3478

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2340.1365
$ws.Range("I134").Value = 1967.4667
$ws.Range("K134").Value = 5902.4001
$ws.Range("M134").Value = -3367.4001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1976.8182
$ws.Range("J16").Value = 2334.75
$ws.Range("L16").Value = 2334.75
$ws.Range("N16").Value = -2908.75
$ws.Range("H31").Value = 2554937.5
$ws.Range("J31").Value = 5686923.5
$ws.Range("L31").Value = 5686923.5
$ws.Range("N31").Value = -5687513.5
$ws.Range("H34").Value = 2554937.5
$ws.Range("J34").Value = 5686923.5
$ws.Range("L34").Value = 5686923.5
$ws.Range("N34").Value = -5687327.5
$ws.Range("H113").Value = 1976.8182
$ws.Range("J113").Value = 2334.75
$ws.Range("L113").Value = 2334.75
$ws.Range("N113").Value = -6674.75
$ws.Range("H122").Value = 419.7
$ws.Range("I122").Value = 384.73334
$ws.Range("K122").Value = 1154.20002
$ws.Range("M122").Value = 1295.79998
$ws.Range("H134").Value = 3845.7646
$ws.Range("I134").Value = 4361.4585
$ws.Range("K134").Value = 13084.3755
$ws.Range("M134").Value = -10549.3755

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1008.4
$ws.Range("I5").Value = 1164
$ws.Range("J5").Value = 775
$ws.Range("K5").Value = 3492
$ws.Range("L5").Value = 2325
$ws.Range("M5").Value = -3380
$ws.Range("N5").Value = -2549
$ws.Range("H6").Value = 77.166664
$ws.Range("I6").Value = 87.59999999999999
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 262.8
$ws.Range("L6").Value = 75
$ws.Range("M6").Value = -149.8
$ws.Range("N6").Value = -301
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()
$ws.Range("H132").Value = 1894.6666
$ws.Range("I132").Value = 1283
$ws.Range("J132").Value = 2200.5
$ws.Range("K132").Value = 11547
$ws.Range("L132").Value = 19804.5
$ws.Range("M132").Value = -9017
$ws.Range("N132").Value = -24864.5
$ws.Range("H134").Value = 3172.8333
$ws.Range("I134").Value = 838.9
$ws.Range("K134").Value = 2516.7
$ws.Range("M134").Value = 2553.3
$ws.Range("H135").Value = 1008.4
$ws.Range("I135").Value = 1164
$ws.Range("J135").Value = 775
$ws.Range("K135").Value = 10476
$ws.Range("L135").Value = 6975
$ws.Range("M135").Value = -7941
$ws.Range("N135").Value = -12045

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 9999
$ws.Range("I21").Value = 9999
$ws.Range("K21").Value = 9999
$ws.Range("M21").Value = -9826
$ws.Range("H30").Value = 9999
$ws.Range("I30").Value = 9999
$ws.Range("K30").Value = 9999
$ws.Range("M30").Value = -9894
$ws.Range("H107").Value = 10321.9
$ws.Range("I107").Value = 3433
$ws.Range("J107").Value = 14914.5
$ws.Range("K107").Value = 3433
$ws.Range("L107").Value = 14914.5
$ws.Range("M107").Value = -1513
$ws.Range("N107").Value = -18754.5
$ws.Range("H132").Value = 1973.3922
$ws.Range("I132").Value = 1468.2609
$ws.Range("K132").Value = 4404.7827
$ws.Range("M132").Value = -1874.7827

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 50000000
$ws.Range("I23").Value = 50000000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 50000000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -49999770
$ws.Range("N23").ClearContents()
$ws.Range("H122").Value = 12012.75
$ws.Range("I122").Value = 11249
$ws.Range("J122").Value = 12267.333
$ws.Range("K122").Value = 33747
$ws.Range("L122").Value = 36801.999
$ws.Range("M122").Value = -31297
$ws.Range("N122").Value = -41701.999
$ws.Range("H132").Value = 4743.4116
$ws.Range("I132").Value = 3913.3125
$ws.Range("J132").Value = 5481.278
$ws.Range("K132").Value = 11739.9375
$ws.Range("L132").Value = 16443.834
$ws.Range("M132").Value = -9209.9375
$ws.Range("N132").Value = -21503.834
$ws.Range("H136").Value = 4226.8184
$ws.Range("J136").Value = 5392.7144
$ws.Range("L136").Value = 16178.1432
$ws.Range("N136").Value = -21278.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 38992
$ws.Range("J140").Value = 38992
$ws.Range("L140").Value = 38992
$ws.Range("N140").Value = -49352
